$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.147.57"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "1.792.33"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'227.62"
$ws.Range("E5").Value = "  +1.87%  "

# Row 6
$ws.Range("E6").Value = "  -0.63%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").Value = "'32.34"
$ws.Range("E8").Value = "  +0.25%  "

# Row 9
$ws.Range("E9").Value = "  +4.63%  "

# Row 10
$ws.Range("D10").Value = "'0.0688"
$ws.Range("E10").Value = "  -2.61%  "

# Row 11
$ws.Range("D11").Value = "'0.0941"
$ws.Range("E11").Value = "  +1.23%  "

# Row 12
$ws.Range("D12").Value = "2.049.51"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("D13").Value = "'11.60"
$ws.Range("E13").Value = "  +6.99%  "

# Row 14
$ws.Range("D14").Value = "1.793.87"
$ws.Range("E14").Value = "  +0.15%  "

# Row 15
$ws.Range("E15").Value = "  +0.27%  "

# Row 16
$ws.Range("D16").Value = "34.114.88"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17
$ws.Range("D17").Value = "'4.20"
$ws.Range("E17").Value = "  +1.32%  "

# Row 18
$ws.Range("D18").Value = "'68.17"
$ws.Range("E18").Value = "  +0.49%  "

# Row 19
$ws.Range("D19").Value = "'244.17"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0778"
$ws.Range("E20").Value = "  -0.40%  "

# Row 21
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").Value = "'10.85"
$ws.Range("E22").Value = "  +1.81%  "

# Row 23
$ws.Range("E23").Value = "  +1.44%  "

# Row 24
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  -1.89%  "

# Row 25
$ws.Range("D25").Value = "'162.06"
$ws.Range("E25").Value = "  +2.19%  "

# Row 26
$ws.Range("D26").Value = "'7.19"
$ws.Range("E26").Value = "  +2.79%  "

# Row 27
$ws.Range("D27").Value = "'16.30"
$ws.Range("E27").Value = "  +0.36%  "

# Row 28
$ws.Range("E28").Value = "  +1.52%  "

# Row 29
$ws.Range("E29").Value = "  -0.11%  "

# Row 30
$ws.Range("E30").Value = "  +2.73%  "

# Row 31
$ws.Range("E31").Value = "  +0.12%  "

# Row 32
$ws.Range("D32").Value = "'3.68"
$ws.Range("E32").Value = "  +0.69%  "

# Row 33
$ws.Range("E33").Value = "  +4.52%  "

# Row 34
$ws.Range("E34").Value = "  +2.40%  "

# Row 35
$ws.Range("D35").Value = "1.411.11"
$ws.Range("E35").Value = "  +1.95%  "

# Row 36
$ws.Range("D36").Value = "'0.657"
$ws.Range("E36").Value = "  +2.05%  "

# Row 37
$ws.Range("E37").Value = "  +0.03%  "

# Row 38
$ws.Range("D38").Value = "'0.0190"
$ws.Range("E38").Value = "  +2.92%  "

# Row 39
$ws.Range("D39").Value = "'2.35"
$ws.Range("E39").Value = "  +8.81%  "

# Row 40
$ws.Range("D40").Value = "'80.48"
$ws.Range("E40").Value = "  +1.53%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.927"
$ws.Range("E41").Value = "  +1.67%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.36"
$ws.Range("E42").Value = "  +0.40%  "

# Row 43
$ws.Range("D43").Value = "'2.69"
$ws.Range("E43").Value = "  -0.53%  "

# Row 44
$ws.Range("D44").Value = "'13.39"
$ws.Range("E44").Value = "  +12.16%  "

# Row 45
$ws.Range("E45").Value = "  +2.81%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'6.08"
$ws.Range("E46").Value = "  +4.29%  "

# Row 47
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0509"
$ws.Range("E47").Value = "  +2.05%  "

# Row 48
$ws.Range("E48").Value = "  +2.41%  "

# Row 49
$ws.Range("D49").Value = "'107.71"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50
$ws.Range("D50").Value = "1.950.99"
$ws.Range("E50").Value = "  -0.33%  "

# Row 51
$ws.Range("E51").Value = "  -0.14%  "
